$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "2007年" data row (row 2). All rows below it shift up by one,
# so "2010年"/"2012年"/"2015年"/"2017年" move from rows 3-6 into rows 2-5.
$ws.Rows.Item(2).Delete()
